$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "2 of 3"

$ws1.Range("B7").Value = "row with a boolean formula"
$ws1.Range("B8").Value = "row with formulas"
$ws1.Range("C8").Formula = "=AVERAGE(2,4,8,16,32,64)"
$ws1.Range("D8").Formula = '=CONCATENATE("string ","cat")'

foreach ($r in @(1,2,3,5,6,7,8,9,10)) {
  $ws1.Rows.Item($r).RowHeight = 12.8
}
foreach ($r in @(1,2,3)) {
  $ws2.Rows.Item($r).RowHeight = 12.8
}

$ws3.Range("A1").Value = "a sheet"
$ws3.Range("B1").Value = "with"
$ws3.Range("C1").Value = 4
$ws3.Range("D1").Value = "columns"
$ws3.Range("A1:D1").Font.Bold = $true

$ws3.Range("A2").Value = "row 1,col 1"
$ws3.Range("B2").Value = "row 1,col 2"
$ws3.Range("C2").Value = "row 1,col 3"
$ws3.Range("D2").Value = "row 1,col 4"

$ws3.Range("A3").Value = "row 2,col 1"
$ws3.Range("B3").Value = "row 2,col 2"
$ws3.Range("C3").Value = "row 2,col 3"
$ws3.Range("D3").Value = "row 2,col 4"

[void]$ws1.Range("B9").Select()
[void]$ws3.Range("C16").Select()
[void]$ws3.Activate()
